$d = $word.ActiveDocument

# The document has several " <en-dash> " separators (after "Course Name",
# "Collage Name", and "Batch number"). Anchor on the "Batch number" label
# and search forward (no wraparound) from there so we land on the correct
# en-dash instead of the first one in the document.
$label = $d.Content.Duplicate
$foundLabel = $label.Find.Execute("Batch number", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$labelEnd = $label.End

$dash = $d.Content.Duplicate
$dash.Start = $labelEnd
$foundDash = $dash.Find.Execute("–", $true, $false, $false, $false, $false, $true, 0, $false, "", 0)
$dashStart = $dash.Start
$dashEnd = $dash.End

# Insert the new text right after the en-dash run.
$ins = $dash.Duplicate
$ins.Collapse(0)  # wdCollapseEnd
$ins.InsertAfter(" 01")

# InsertAfter merges the new text into whichever neighbouring runs share
# its formatting (here every run on the line is sz=36pt/"en-US"), which
# would also fuse the untouched " " / "-" / " " runs around it into one
# big run. Re-split the four logical segments back into their own runs by
# touching each with a differing-then-restored font size (matching the
# diff's sz=72/szCs=72 half-points == Font.Size 36pt).
$segSpaceBefore = $d.Range($labelEnd, $dashStart)     # " "  (untouched)
$segDash        = $d.Range($dashStart, $dashEnd)      # "-"  (untouched)
$segNew         = $d.Range($dashEnd, $dashEnd + 3)    # " 01" (new run)
$segSpaceAfter  = $d.Range($dashEnd + 3, $dashEnd + 4) # " "  (untouched)

foreach ($seg in @($segSpaceBefore, $segDash, $segNew, $segSpaceAfter)) {
    $seg.Font.Size = 11
}
foreach ($seg in @($segSpaceBefore, $segDash, $segNew, $segSpaceAfter)) {
    $seg.Font.Size = 36
}
